# Add example: empty object
# - Adds a new "field7" / "ref" / "seasonEmptyField" row to the Lvl0 sheet
#   (mirrors the existing field-description rows) and a corresponding
#   "seasonEmptyField" entry on the Lvl1 sheet, representing an empty object
#   example value.

$wb = $excel.ActiveWorkbook

# --- Lvl0 ("B1:F26") -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Lvl0")

$ws1.Range("D19").Value = "field7"
$ws1.Range("E19").Value = "ref"
$ws1.Range("F19").Value = "seasonEmptyField"

# --- Lvl1 ("B1:F21" -> "B1:F23") --------------------------------------
$ws2 = $wb.Worksheets.Item("Lvl1")

$ws2.Range("C23").Value = "seasonEmptyField"

# Update the selections last, finishing on the sheet that was already the
# active one (Lvl1) so the workbook's active-tab state is preserved.
$ws1.Range("D20").Select()
$ws2.Range("C25").Select()
